# Insert 3 new weekly report rows before the existing row 215, shifting
# rows 215:317 down to 218:320 (dimension grows from A1:T317 to A1:T320).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A215:A217").EntireRow.Insert()

# Shared (static) column values for this Comercializadora del Agro de
# Limarí / Palta block - identical across the whole sheet.
$colA = 2
$colB = "Comercializadora del Agro de Limarí"
$colC = "Coquimbo"
$colE = 4
$colF = "Fruta"
$colG = 100106
$colH = "Oleaginosos"
$colI = 100106002
$colJ = "Palta"
$colR = "Provincia de Limarí"
$colT = 1

# New data rows (date 2021-12-23 -> serial 44553), variety Hass, quality
# Especial / Primera / Segunda, box of 17 kilos.
$newRows = @(
    @{ Row = 215; L = "Especial"; M = 560; N = 2400; O = 2500; P = 2450 },
    @{ Row = 216; L = "Primera";  M = 400; N = 2100; O = 2200; P = 2150 },
    @{ Row = 217; L = "Segunda";  M = 300; N = 1700; O = 1800; P = 1750 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $colA
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = 44553
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
    $ws.Cells.Item($row, 10).Value = $colJ
    $ws.Cells.Item($row, 11).Value = "Hass"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/kilo (en caja de 17 kilos)"
    $ws.Cells.Item($row, 18).Value = $colR
    $ws.Cells.Item($row, 19).Value = $r.P
    $ws.Cells.Item($row, 20).Value = $colT
}
